# Update the test fixture workbook:
#  - rename the administrative-division header labels on the "School"
#    sheet (County/Sub-County/Ward -> Province/District/Subdistrict)
#  - move the selection on "School" to G1
#  - make "School" the active sheet/tab (it was "Health Care Facilities")
$wb = $excel.ActiveWorkbook

$wsSchool = $wb.Worksheets.Item("School")

$wsSchool.Range("D1").Value = "Province"
$wsSchool.Range("E1").Value = "District"
$wsSchool.Range("F1").Value = "Subdistrict"

$wsSchool.Range("G1").Select()
$wsSchool.Activate()
